$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 56 ("Freemarker" section continues): REST entry gets a parenthetical
# and its own numbered video links; date shifts from 44314 -> 44317.
# ---------------------------------------------------------------------------
$ws.Range("C56").Value = "REST (Json, XML)"
$ws.Range("E56").Value = 44317
$ws.Range("F56").Value = "1. https://youtu.be/wNYuN-5TcCk `n2. https://youtu.be/2VSN0CZhTJE"

# Copy the wrapped-hyperlink-style formatting already used elsewhere (F54)
# onto F56 so it renders the same way as the other multi-link cells.
$ws.Range("F54").Copy() | Out-Null
$ws.Range("F56").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Old row 57 ("Web project in Java EE", #8.2) is replaced entirely by the
# new lesson rows below, so wipe it completely (values + formatting).
$ws.Rows.Item(57).Clear() | Out-Null

# ---------------------------------------------------------------------------
# New rows 58-65: "Web project in Java EE" is split into Part #1 / Part #2,
# each with its own numbered sub-lessons (lesson #41 and #42).
# ---------------------------------------------------------------------------
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = "Web project in Java EE `u{2013} Part #1"
$ws.Range("D58").Value = 2
$ws.Range("E58").Value = 44319
$ws.Range("F58").Value = "1. https://youtu.be/q1a3CS49zHE `n2. https://youtu.be/6G_fIVU4VKM `n3. https://youtu.be/jhAcrXZHab0 "

$ws.Range("B59").Value = 2
$ws.Range("C59").Value = "Web project in Java EE `u{2013} Part #2"
$ws.Range("D59").Value = 2

$ws.Range("B60").Value = 3
$ws.Range("D60").Value = 2

$ws.Range("B61").Value = 4
$ws.Range("D61").Value = 2

$ws.Range("B62").Value = 5
$ws.Range("D62").Value = 2

$ws.Range("B63").Value = 6
$ws.Range("D63").Value = 2

$ws.Range("B64").Value = 7
$ws.Range("D64").Value = 2

$ws.Range("B65").Value = 8
$ws.Range("D65").Value = 2

# ---------------------------------------------------------------------------
# Formatting for the new rows - reuse existing style "donor" cells so the
# new rows look consistent with the rest of the table.
# ---------------------------------------------------------------------------
# Column B: first sub-item of a lesson group is bold (style of B4), the rest
# use the plain numbered style (style of D4).
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B58").PasteSpecial(-4122) | Out-Null
$ws.Range("B65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D4").Copy() | Out-Null
$ws.Range("B59:B64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column C: lesson-name style (style of C10).
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C58:C65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column D: hours style (style of D4).
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D58:D65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column E: date style (style of E6).
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E58:E65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column F: link-cell styles, matching the pattern used for the other
# youtube-link cells in this block.
$ws.Range("F44").Copy() | Out-Null
$ws.Range("F58").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F50").Copy() | Out-Null
$ws.Range("F59:F61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F44").Copy() | Out-Null
$ws.Range("F62:F63").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F26").Copy() | Out-Null
$ws.Range("F64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F47").Copy() | Out-Null
$ws.Range("F65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row heights for the new content rows (auto heights observed in the source
# authoring tool once the cells hold their final wrapped text).
$ws.Rows.Item(56).RowHeight = 28.45
$ws.Rows.Item(58).RowHeight = 41.75
$ws.Rows.Item(59).RowHeight = 13.8
$ws.Rows.Item(60).RowHeight = 13.8
$ws.Rows.Item(61).RowHeight = 13.8
$ws.Rows.Item(62).RowHeight = 13.8
$ws.Rows.Item(63).RowHeight = 13.8
$ws.Rows.Item(64).RowHeight = 13.8
$ws.Rows.Item(65).RowHeight = 13.8

# ---------------------------------------------------------------------------
# The old single-video hyperlink on F55 no longer applies now that its text
# holds three numbered links, so drop it.
# ---------------------------------------------------------------------------
$ws.Range("F55").Hyperlinks.Delete()

# Move the view to roughly where the editor left it.
$ws.Range("F58").Select()
